# Amplifier.xlsx -- "Formula for lm35 with ampifier"
#
# Adds a small offset note in A4 (new shared string), restyles the G
# column (Ua via formula, in mV) to use the "0 mV" number format instead
# of "0.000 V", widens column A to fit the new label, and moves the
# saved selection back to A4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- A4: new label " (909mV offset)" in a small (9pt), non-bold font
#     with the same bottom border used by the rest of the header row.
$ws.Range("A4").Value = " (909mV offset)"
$ws.Range("A4").Font.Bold = $false
$ws.Range("A4").Font.Size = 9
$ws.Range("A4").Borders.Item(9).LineStyle = 1

# --- Column A: widen slightly so the new label fits.
$ws.Columns("A").ColumnWidth = 11.6

# --- G5:G12 ("Ua" computed via formula, in mV): switch number format
#     from volts ("0.000 V") to millivolts ("0 mV"), matching the other
#     mV columns (C, E).
$ws.Range("G5:G12").NumberFormat = '0\ "mV"'

# --- Restore the saved selection/view to A4 (was E38 with the sheet
#     scrolled down before).
$ws.Activate()
$ws.Range("A4").Select()
